$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are written as exact text (avoid numeric/locale coercion
# of values like "1.002" or trailing-zero loss like "2.650" -> 2.65).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.639.07"
$ws.Range("E2").Value = "  +1.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.83"
$ws.Range("E3").Value = "  +1.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.78"

$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4823"
$ws.Range("E7").Value = "  +0.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2866"
$ws.Range("E8").Value = "  +2.11%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06548"
$ws.Range("E9").Value = "  +0.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.917.80"
$ws.Range("E10").Value = "  +3.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07473"
$ws.Range("E11").Value = "  +1.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.68"
$ws.Range("E12").Value = "  +2.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.103"
$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.08"
$ws.Range("E14").Value = "  +0.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6665"
$ws.Range("E15").Value = "  +2.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.639.89"
$ws.Range("E16").Value = "  +1.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.29"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007574"
$ws.Range("E19").Value = "  -0.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "232.11"
$ws.Range("E20").Value = "  +3.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.157.62"
$ws.Range("E21").Value = "  +2.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.282"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("E24").Value = "  +2.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.71"
$ws.Range("E25").Value = "  +3.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.347"
$ws.Range("E26").Value = "  +0.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.77"
$ws.Range("E27").Value = "  +1.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.964"
$ws.Range("E28").Value = "  +2.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.401"
$ws.Range("E29").Value = "  -2.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1019"
$ws.Range("E30").Value = "  +10.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.343"
$ws.Range("E31").Value = "  +2.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.025"
$ws.Range("E32").Value = "  +1.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05099"
$ws.Range("E33").Value = "  +1.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.213"
$ws.Range("E34").Value = "  +5.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7585"
$ws.Range("E35").Value = "  +2.46%  "

$ws.Range("E36").Value = "  +0.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01876"
$ws.Range("E37").Value = "  +2.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.650"
$ws.Range("E38").Value = "  +1.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9204"
$ws.Range("E39").Value = "  +1.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.069"
$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.09"
$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4302"
$ws.Range("E42").Value = "  +0.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.004"
$ws.Range("E43").Value = "  +0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.655"
$ws.Range("E44").Value = "  -5.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.430"
$ws.Range("E45").Value = "  +0.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.28"
$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1275"
$ws.Range("E47").Value = "  -3.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.487"
$ws.Range("E48").Value = "  -3.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.011"
$ws.Range("E49").Value = "  +3.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.92"
$ws.Range("E50").Value = "  -0.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05668"
$ws.Range("E51").Value = "  +0.03%  "
Write-Output "applied 93 cell updates"
